$d = $word.ActiveDocument

# Locate the two bibliography paragraphs that need to be merged:
#   1) "Zill, Dennis, Álgebra y trigonometría. Ed McGraw Hill"
#   2) "Buriticá Trujillo, B. (2009). Álgebra y trigonometría."
$zillPara = $null
$buriticaPara = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    $text = $para.Range.Text
    if ($text -like "*Zill, Dennis*") {
        $zillPara = $para
    }
    if ($text -like "*Buritic*Trujillo*") {
        $buriticaPara = $para
    }
}

if ($zillPara -eq $null) {
    throw "Could not find the 'Zill, Dennis' bibliography paragraph"
}
if ($buriticaPara -eq $null) {
    throw "Could not find the 'Buriticá Trujillo' bibliography paragraph"
}

# Replace the first paragraph's text with the new reference, split into
# three runs: a plain-text lead-in, an italicised book title, and a
# trailing period.
$prefix = "Varberg, P., & Rigdon, R. P. V. (2000). Cálculo diferencial e integral. "
$italicPart = "Edit. Prince Hall"
$suffix = "."

$zillPara.Range.Text = $prefix + $italicPart + $suffix

$italicStart = $zillPara.Range.Start + $prefix.Length
$italicEnd = $italicStart + $italicPart.Length
$italicRange = $d.Range($italicStart, $italicEnd)
$italicRange.Italic = 1

# Remove the now-redundant "Buriticá Trujillo" paragraph entirely
# (its content has been folded into the reference above).
$buriticaPara.Range.Delete()
